# Insert a new data row at row 27 (shifts existing rows 27..59 down to 28..60)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(27).Insert()

# Populate the newly inserted row 27 with the new weekly record
$ws.Range("A27").Value = 9
$ws.Range("B27").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C27").Value = "Metropolitana"
$ws.Range("D27").Value = 44539
$ws.Range("E27").Value = 13
$ws.Range("F27").Value = "Fruta"
$ws.Range("G27").Value = 100101
$ws.Range("H27").Value = "Berries"
$ws.Range("I27").Value = 100101004
$ws.Range("J27").Value = "Frambuesa"
$ws.Range("K27").Value = "Sin especificar"
$ws.Range("L27").Value = "Primera"
$ws.Range("M27").Value = 450
$ws.Range("N27").Value = 6000
$ws.Range("O27").Value = 6000
$ws.Range("P27").Value = 6000
$ws.Range("Q27").Value = "$/bandeja 2 kilos"
$ws.Range("R27").Value = "Provincia de Curicó"
$ws.Range("S27").Value = 3000
$ws.Range("T27").Value = 2
